$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split a run at a given absolute character offset by inserting a
# zero-length bookmark there and immediately deleting it. The bookmark
# insertion forces the text-run to be broken in two at that point (even when
# the formatting on both sides is identical), and once the bookmark itself
# is deleted the break persists, giving us two separate <w:r> elements.
# ---------------------------------------------------------------------------
$script:splitCounter = 0
function Split-RunAt([int]$pos) {
    $script:splitCounter += 1
    $name = "tmpSplitMarker$($script:splitCounter)"
    $d.Bookmarks.Add($name, $d.Range($pos, $pos)) | Out-Null
    $d.Bookmarks($name).Delete()
}

# ---------------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that originally sat between
# the "0-20" run and the " min)" run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: "Questions to the presenter (~10 min)" -> "Q" + "&A with " +
# "the presenter (~10 min)"  (i.e. "Q&A with the presenter (~10 min)")
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Questions to the presenter (~10 min)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$rng.Text = "Q&A with the presenter (~10 min)"
Split-RunAt ($start + 1)
Split-RunAt ($start + 9)

# ---------------------------------------------------------------------------
# Change 3: "...(~2" + "0 min)" -> "...(~2" + "0 min" + ", random assignment
# to breakout room)"
# (use the unique phrase "(~20 min)" to anchor the search, then only touch
# the trailing "0 min)" portion of that match)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Thoughts on topic (~20 min)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchEnd = $rng.End
$start = $matchEnd - 6   # start of the "0 min)" run
$target = $d.Range($start, $matchEnd)
$target.Text = "0 min, random assignment to breakout room)"
Split-RunAt ($start + 5)

# ---------------------------------------------------------------------------
# Change 4: "pdates on projects/concerns related to work (~30 min)" ->
# "pdates on projects/concerns related to work (~30 min" +
# ", assigned to break out room based on experience/challenge level?" +
# [bookmark _GoBack] + ")"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("pdates on projects/concerns related to work (~30 min)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$rng.Text = "pdates on projects/concerns related to work (~30 min, assigned to break out room based on experience/challenge level?)"
$p1 = $start + 52
$p2 = $start + 52 + 65
Split-RunAt $p1
Split-RunAt $p2
$d.Bookmarks.Add("_GoBack", $d.Range($p2, $p2)) | Out-Null
